# Updates the cryptos price/volume table to the latest scraped snapshot.
# Mirrors the GitHub Actions "Updated cryptos list" commit: refreshed
# Price (D) / Volume(1h) (E) values for each row, and for rows 20-21 the
# Uniswap / InternetComputer(DFINITY) entries swapped rank order (so the
# full row -- Coin, Link, Price, Volume -- moves together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (cell, new text) pairs exactly as they should appear in the sheet.
# Kept as literal text (not auto-converted to numbers/formulas) because
# the source data uses "." as a thousands separator (e.g. "51.912.77")
# and fixed-width percent strings with padding spaces (e.g. "  +4.84%  ").
$updates = @(
    @{ Cell='D2'; Value='51.912.77' }
    @{ Cell='E2'; Value='  +4.84%  ' }
    @{ Cell='D3'; Value='2.771.45' }
    @{ Cell='E3'; Value='  +4.94%  ' }
    @{ Cell='E4'; Value='  +0.05%  ' }
    @{ Cell='D5'; Value='339.58' }
    @{ Cell='E5'; Value='  +4.10%  ' }
    @{ Cell='D6'; Value='115.34' }
    @{ Cell='E6'; Value='  +2.18%  ' }
    @{ Cell='D7'; Value='0.544' }
    @{ Cell='E7'; Value='  +4.21%  ' }
    @{ Cell='D8'; Value='0.999' }
    @{ Cell='E8'; Value='  -0.04%  ' }
    @{ Cell='E9'; Value='  +4.48%  ' }
    @{ Cell='D10'; Value='41.63' }
    @{ Cell='E10'; Value='  +4.92%  ' }
    @{ Cell='D11'; Value='0.0858' }
    @{ Cell='E11'; Value='  +5.36%  ' }
    @{ Cell='D12'; Value='19.94' }
    @{ Cell='E12'; Value='  -0.41%  ' }
    @{ Cell='E13'; Value='  +1.89%  ' }
    @{ Cell='D14'; Value='7.57' }
    @{ Cell='E14'; Value='  -0.06%  ' }
    @{ Cell='D15'; Value='3.206.14' }
    @{ Cell='E15'; Value='  +4.98%  ' }
    @{ Cell='D16'; Value='2.777.01' }
    @{ Cell='E16'; Value='  +5.92%  ' }
    @{ Cell='D17'; Value='51.759.84' }
    @{ Cell='E17'; Value='  +4.57%  ' }
    @{ Cell='D18'; Value='0.875' }
    @{ Cell='E18'; Value='  +1.91%  ' }
    @{ Cell='D19'; Value='3.18' }
    @{ Cell='E19'; Value='  +9.34%  ' }
    @{ Cell='B20'; Value='Uniswap' }
    @{ Cell='C20'; Value='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' }
    @{ Cell='D20'; Value='6.94' }
    @{ Cell='E20'; Value='  +4.14%  ' }
    @{ Cell='B21'; Value='InternetComputer(DFINITY)' }
    @{ Cell='C21'; Value='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell='D21'; Value='13.17' }
    @{ Cell='E21'; Value='  -1.12%  ' }
    @{ Cell='D22'; Value='0.0₃0975' }
    @{ Cell='E22'; Value='  +2.67%  ' }
    @{ Cell='D23'; Value='276.17' }
    @{ Cell='E23'; Value='  +3.05%  ' }
    @{ Cell='D24'; Value='69.70' }
    @{ Cell='E24'; Value='  +0.78%  ' }
    @{ Cell='E25'; Value='  +6.94%  ' }
    @{ Cell='E26'; Value='  +2.10%  ' }
    @{ Cell='D27'; Value='1.00' }
    @{ Cell='E27'; Value='  +0.00%  ' }
    @{ Cell='D28'; Value='10.14' }
    @{ Cell='E28'; Value='  +0.16%  ' }
    @{ Cell='E29'; Value='  -0.64%  ' }
    @{ Cell='D30'; Value='0.140' }
    @{ Cell='E30'; Value='  +2.28%  ' }
    @{ Cell='D31'; Value='34.46' }
    @{ Cell='E31'; Value='  -0.56%  ' }
    @{ Cell='D32'; Value='50.10' }
    @{ Cell='E32'; Value='  +0.91%  ' }
    @{ Cell='D33'; Value='5.69' }
    @{ Cell='E33'; Value='  +4.01%  ' }
    @{ Cell='D34'; Value='0.0820' }
    @{ Cell='E34'; Value='  -0.02%  ' }
    @{ Cell='E35'; Value='  +0.02%  ' }
    @{ Cell='E36'; Value='  +2.98%  ' }
    @{ Cell='D37'; Value='18.84' }
    @{ Cell='E37'; Value='  -1.44%  ' }
    @{ Cell='D38'; Value='4.92' }
    @{ Cell='E38'; Value='  -0.33%  ' }
    @{ Cell='D39'; Value='3.20' }
    @{ Cell='E39'; Value='  +2.35%  ' }
    @{ Cell='D40'; Value='0.0374' }
    @{ Cell='E40'; Value='  +10.34%  ' }
    @{ Cell='D41'; Value='2.65' }
    @{ Cell='E41'; Value='  +25.71%  ' }
    @{ Cell='D42'; Value='2.33' }
    @{ Cell='E42'; Value='  +1.32%  ' }
    @{ Cell='E43'; Value='  +3.19%  ' }
    @{ Cell='D44'; Value='125.82' }
    @{ Cell='E44'; Value='  -2.79%  ' }
    @{ Cell='D45'; Value='23.14' }
    @{ Cell='E45'; Value='  -1.65%  ' }
    @{ Cell='D46'; Value='2.069.76' }
    @{ Cell='E46'; Value='  +0.39%  ' }
    @{ Cell='E47'; Value='  -0.26%  ' }
    @{ Cell='E48'; Value='  +1.95%  ' }
    @{ Cell='D49'; Value='5.50' }
    @{ Cell='E49'; Value='  +4.86%  ' }
    @{ Cell='D50'; Value='8.82' }
    @{ Cell='E50'; Value='  -0.86%  ' }
    @{ Cell='D51'; Value='59.03' }
    @{ Cell='E51'; Value='  +0.31%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force Text format before assigning so Excel's COM layer doesn't
    # reinterpret numeric-looking strings (e.g. "339.58", "1.00") as
    # numbers, then restore the default "Normal" style so no stray
    # per-cell formatting is left behind (matches the original workbook,
    # where these cells carry no explicit style).
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
